# Update the "Metadata" sheet: URL, Version, Date, Publisher moved from the
# old ibm.com / Alvearie naming to the new linuxforhealth.org / LinuxForHealth
# naming, plus the version bump and new publish date.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/legally-documented-sex"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet: the three nested extension "Type(s)" cells and
# the "Fixed Value" URL cell need the same ibm.com -> linuxforhealth.org
# rename, and the top-level Extension row's Constraint(s) cell is cleared.
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("AI2").Value = ""
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-sex}`n"
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-source}`n"
$elements.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-system}`n"
$elements.Range("Q8").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/legally-documented-sex"
